# ErgoLux Swedish (sv-SE) translation workbook update
# - Add a new translation row for the "strWindowPos" key (settings -> User
#   interface tab): "Remember window position and size on startup", with a
#   matching Comment.
# - Give the existing "strChkDlgPath" row (row 25) the same Comment, since it
#   also lives on the settings form's "User interface" tab.
# - Keep the "Tabla13" table range / AutoFilter in sync with the newly
#   inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Comment for the pre-existing "strChkDlgPath" entry.
$ws.Range("D25").Value = 'In "settings" form, tab "User interface"'

# Insert a brand-new row for "strWindowPos" right after row 31 (before the
# old row 32 / "strDlgReset"), shifting everything below it down by one.
$ws.Rows.Item(32).Insert()

# Grow the table so the new row becomes part of "Tabla13" (updates both the
# table ref and its AutoFilter range).
$lo.Resize($ws.Range("B2:F204"))

$ws.Range("B32").Value = "localization\strings"
$ws.Range("C32").Value = "strWindowPos"
$ws.Range("D32").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E32").Value = "Remember window position and size on startup"

# The Comment column got a little wider to fit the new text.
$ws.Columns.Item(4).ColumnWidth = 34.8
